$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.302.17"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "3.723.58"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'612.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.31%  "
$ws.Range("D6").Value = "'191.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.48%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.728"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'60.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.77%  "
$ws.Range("D11").Value = "'0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("D13").Value = "'10.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "4.319.11"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "3.725.56"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'19.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "69.121.79"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'412.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'4.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'89.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'3.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "'12.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'3.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'6.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "'9.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'33.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "'7.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Value = "'12.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").Value = "'45.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("D35").Value = "'635.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'66.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "'0.418"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("D38").Value = "0.0₃0825"
$ws.Range("E38").Value = "  -11.01%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "'0.0448"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").Value = "2.882.15"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("D47").Value = "'9.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'142.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'3.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("E51").Value = "  -0.43%  "
